$wb = $excel.ActiveWorkbook

# Update 2024 (and a couple 2023 correction) violent-crime figures for 2024-11-21 data refresh.
# Each worksheet is addressed by its tab name; cells are set directly by A1-style reference.

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 7187
$ws.Range("K3").Value = 7451
$ws.Range("J4").Value = 1845
$ws.Range("K4").Value = 1552
$ws.Range("K5").Value = 526
$ws.Range("K6").Value = 8231
$ws.Range("J7").Value = 29313
$ws.Range("K7").Value = 24947

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 449
$ws.Range("K6").Value = 543
$ws.Range("K7").Value = 1620

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K6").Value = 337
$ws.Range("K7").Value = 1063

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 141
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 410

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 275
$ws.Range("K7").Value = 836

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 140
$ws.Range("K4").Value = 25
$ws.Range("K6").Value = 220
$ws.Range("K7").Value = 591

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 217
$ws.Range("K7").Value = 752
$ws.Range("K8").Value = 1620
$ws.Range("K12").Value = 44
$ws.Range("K18").Value = 165
$ws.Range("K19").Value = 726
$ws.Range("K20").Value = 614
$ws.Range("K21").Value = 87
$ws.Range("K23").Value = 255
$ws.Range("K25").Value = 115
$ws.Range("K29").Value = 1372
$ws.Range("K31").Value = 288
$ws.Range("K33").Value = 1063
$ws.Range("K37").Value = 836
$ws.Range("K41").Value = 171
$ws.Range("K42").Value = 921
$ws.Range("K43").Value = 207
$ws.Range("K44").Value = 205
$ws.Range("K47").Value = 167
$ws.Range("K48").Value = 319
$ws.Range("K51").Value = 315
$ws.Range("K52").Value = 648
$ws.Range("K54").Value = 486
$ws.Range("K55").Value = 269
$ws.Range("K60").Value = 143
$ws.Range("K63").Value = 66
$ws.Range("K64").Value = 147
$ws.Range("K65").Value = 591
$ws.Range("K67").Value = 981
$ws.Range("K70").Value = 44
$ws.Range("K71").Value = 75
$ws.Range("K72").Value = 122
$ws.Range("K73").Value = 224
$ws.Range("K76").Value = 345
$ws.Range("K77").Value = 166
$ws.Range("K78").Value = 301
$ws.Range("K79").Value = 611
$ws.Range("K85").Value = 1145
$ws.Range("K86").Value = 152
$ws.Range("K88").Value = 266
$ws.Range("K89").Value = 375
$ws.Range("K90").Value = 243
$ws.Range("K95").Value = 410
$ws.Range("J96").Value = 329
$ws.Range("K96").Value = 267
$ws.Range("J101").Value = 29313
$ws.Range("K101").Value = 24947

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("K2").Value = 87
$ws.Range("K3").Value = 69
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 275
$ws.Range("K3").Value = 352
$ws.Range("K6").Value = 275
$ws.Range("K7").Value = 981

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 486

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 389
$ws.Range("K3").Value = 489
$ws.Range("K7").Value = 1372

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 43
$ws.Range("K7").Value = 319

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 219
$ws.Range("K6").Value = 239
$ws.Range("K7").Value = 726

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K2").Value = 55
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 205

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 76
$ws.Range("K3").Value = 65
$ws.Range("K4").Value = 26
$ws.Range("K7").Value = 345

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 171

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K2").Value = 251
$ws.Range("K6").Value = 341
$ws.Range("K7").Value = 921

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K5").Value = 7
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 255

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("J4").Value = 21
$ws.Range("K4").Value = 16
$ws.Range("J7").Value = 329
$ws.Range("K7").Value = 267

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 194
$ws.Range("K7").Value = 611

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 211
$ws.Range("K4").Value = 29
$ws.Range("K7").Value = 614

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 165

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K3").Value = 241
$ws.Range("K7").Value = 752

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K3").Value = 40
$ws.Range("K7").Value = 115

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 50
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K3").Value = 116
$ws.Range("K6").Value = 155

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 76
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K2").Value = 65
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 84
$ws.Range("K7").Value = 266

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 117
$ws.Range("K7").Value = 375

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 63
$ws.Range("K7").Value = 152

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 89
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 84
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 143

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 42
$ws.Range("K6").Value = 78
$ws.Range("K7").Value = 207

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 375
$ws.Range("K3").Value = 398
$ws.Range("K6").Value = 283
$ws.Range("K7").Value = 1145

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K3").Value = 27
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 29
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K4").Value = 36
$ws.Range("K7").Value = 648

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 44
